$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(2, 17, 15),
    @(3, 37, 28),
    @(4, 78, 56),
    @(5, 67, 59),
    @(6, 103, 74),
    @(7, 66, 44),
    @(8, 43, 28),
    @(9, 49, 38),
    @(10, 26, 19),
    @(11, 56, 44),
    @(12, 31, 25),
    @(13, 49, 35),
    @(14, 87, 70),
    @(15, 53, 40),
    @(16, 69, 53),
    @(17, 28, 23),
    @(18, 42, 30),
    @(19, 37, 27),
    @(20, 29, 22),
    @(21, 58, 43),
    @(22, 30, 20),
    @(23, 21, 13),
    @(24, 61, 48),
    @(25, 32, 24),
    @(26, 32, 22),
    @(27, 54, 40),
    @(28, 52, 38),
    @(29, 61, 40),
    @(30, 45, 34),
    @(31, 54, 40),
    @(32, 41, 30),
    @(33, 45, 31),
    @(34, 57, 38),
    @(35, 52, 40),
    @(36, 38, 21),
    @(37, 51, 32),
    @(38, 60, 44),
    @(39, 58, 49),
    @(40, 73, 62),
    @(41, 50, 38),
    @(42, 69, 49),
    @(43, 57, 44),
    @(44, 61, 49),
    @(45, 33, 18),
    @(46, 63, 44),
    @(47, 46, 36),
    @(48, 46, 27),
    @(49, 54, 44),
    @(50, 34, 25),
    @(51, 36, 26),
    @(52, 44, 33),
    @(53, 81, 67),
    @(54, 49, 25),
    @(55, 63, 41),
    @(56, 81, 60),
    @(57, 63, 35),
    @(58, 62, 49),
    @(59, 59, 45),
    @(60, 51, 31),
    @(61, 30, 22),
    @(62, 54, 37),
    @(63, 100, 88),
    @(64, 29, 20),
    @(65, 28, 23),
    @(66, 59, 40),
    @(67, 73, 57),
    @(68, 45, 30),
    @(69, 93, 72),
    @(70, 46, 35),
    @(71, 67, 58),
    @(72, 38, 17),
    @(73, 55, 35),
    @(74, 106, 88),
    @(75, 108, 104),
    @(76, 26, 13),
    @(77, 127, 124),
    @(78, 73, 53),
    @(79, 99, 98),
    @(80, 107, 87),
    @(81, 71, 46),
    @(82, 25, 14),
    @(83, 64, 52),
    @(84, 180, 111),
    @(85, 0, 0),
    @(86, 11, 8),
    @(87, 5, 4),
    @(88, 0, 0),
    @(89, 0, 0),
    @(90, 3, 1),
    @(91, 0, 0),
    @(92, 211, 134),
    @(93, 4977, 3696)
)

foreach ($item in $values) {
    $r = $item[0]
    $cval = $item[1]
    $dval = $item[2]
    $ws.Cells.Item($r, 3).Value = $cval
    $ws.Cells.Item($r, 4).Value = $dval
}

$wb.Save()